# Add total deaths plot - update underlying data tables that feed the chart.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Tüm" (index 2): append a new day of data (row 13).
# ---------------------------------------------------------------------------
$wsTum = $wb.Worksheets.Item(2)

$wsTum.Range("A13").Value2 = 43912
$wsTum.Range("B13").Value2 = 947
$wsTum.Range("C13").Formula = "=B13-B12"
$wsTum.Range("D13").Value2 = 21
$wsTum.Range("E13").Formula = "=D13-D12"
$wsTum.Range("F13").Value2 = 12

# Selection moves off the data range after the edit.
$wsTum.Range("M28").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Özet" (index 3): insert a new "Aktif Vaka" column and refresh totals.
# ---------------------------------------------------------------------------
$wsOzet = $wb.Worksheets.Item(3)

# Insert a new column C (shifts old C:F -> D:G, formulas auto-adjust).
$wsOzet.Columns.Item(3).Insert() | Out-Null

$wsOzet.Range("C1").Value2 = "Aktif Vaka"
$wsOzet.Range("C2").Formula = "=A2-B2"

$wsOzet.Range("A2").Value2 = 947
$wsOzet.Range("B2").Value2 = 21

# Restore Özet as the active/selected sheet & its selection.
$wsOzet.Range("A1:E2").Select() | Out-Null
